$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.229.64'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.005.08'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Formula = "'258.92"
$ws.Range("E5").Value = '  +4.64%  '
$ws.Range("D6").Formula = "'0.616"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Formula = "'56.51"
$ws.Range("E8").Value = '  -6.03%  '
$ws.Range("D9").Formula = "'0.380"
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("E10").Value = '  -4.48%  '
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").Value = '2.301.99'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").Formula = "'14.20"
$ws.Range("E13").Value = '  -6.24%  '
$ws.Range("D14").Formula = "'21.65"
$ws.Range("E14").Value = '  -1.87%  '
$ws.Range("E15").Value = '  -7.28%  '
$ws.Range("D16").Formula = "'5.21"
$ws.Range("E16").Value = '  -5.22%  '
$ws.Range("D17").Value = '2.010.84'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '37.284.77'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Formula = "'70.05"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("D20").Value = '0.0₃0832'
$ws.Range("E20").Value = '  -3.53%  '
$ws.Range("D21").Formula = "'233.75"
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("E22").Value = '  -2.66%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").Formula = "'164.74"
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").Formula = "'8.95"
$ws.Range("E27").Value = '  -5.03%  '
$ws.Range("E28").Value = '  -1.04%  '
$ws.Range("D29").Formula = "'0.128"
$ws.Range("E29").Value = '  -7.35%  '
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("E32").Value = '  -4.52%  '
$ws.Range("D33").Formula = "'0.0637"
$ws.Range("E33").Value = '  -5.23%  '
$ws.Range("D34").Formula = "'4.43"
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("E35").Value = '  -5.93%  '
$ws.Range("D36").Formula = "'3.40"
$ws.Range("E36").Value = '  -6.34%  '
$ws.Range("D37").Formula = "'1.81"
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").Formula = "'5.44"
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("D40").Formula = "'3.04"
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("E43").Value = '  -5.84%  '
$ws.Range("D44").Value = '1.435.08'
$ws.Range("E44").Value = '  +4.12%  '
$ws.Range("D45").Formula = "'89.05"
$ws.Range("E45").Value = '  -3.25%  '
$ws.Range("D46").Formula = "'15.58"
$ws.Range("E46").Value = '  -8.60%  '
$ws.Range("E47").Value = '  -3.51%  '
$ws.Range("D48").Formula = "'2.92"
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("D49").Formula = "'6.97"
$ws.Range("E49").Value = '  -6.45%  '
$ws.Range("D50").Value = '2.193.94'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("E51").Value = '  -10.10%  '
